$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2863  # was 2862
$ws.Range("F3").Value = 1160  # was 1159
$ws.Range("G3").Value = 70  # was 60
$ws.Range("F4").Value = 20940  # was 20913
$ws.Range("F6").Value = 2746  # was 2736
$ws.Range("F7").Value = 796  # was 795
$ws.Range("F9").Value = 506  # was 502
$ws.Range("F10").Value = 757  # was 753
$ws.Range("F14").Value = 108  # was 105
$ws.Range("F15").Value = 509  # was 507
$ws.Range("F17").Value = 257  # was 254
$ws.Range("F19").Value = 419  # was 415
$ws.Range("F20").Value = 47  # was 41
$ws.Range("G22").Value = 85  # was 70
$ws.Range("F24").Value = 121  # was 120

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 18  # was 17
$ws.Range("F4").Value = 124  # was 123
$ws.Range("F5").Value = 330  # was 327
$ws.Range("F6").Value = 141  # was 140
$ws.Range("F12").Value = 97  # was 96
$ws.Range("F13").Value = 0  # was 40
$ws.Range("F14").Value = 143  # was 141
$ws.Range("F17").Value = 6  # was 5
$ws.Range("F19").Value = 23  # was 20

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6117  # was 6114
$ws.Range("F4").Value = 675  # was 673
$ws.Range("F5").Value = 1550  # was 1531
$ws.Range("F6").Value = 52  # was 51

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6117  # was 6114
$ws.Range("F4").Value = 675  # was 673
$ws.Range("F5").Value = 1550  # was 1531
$ws.Range("F6").Value = 2863  # was 2862
$ws.Range("F7").Value = 1160  # was 1159
$ws.Range("G7").Value = 70  # was 60
$ws.Range("F8").Value = 20940  # was 20913
$ws.Range("F9").Value = 18  # was 17
$ws.Range("F12").Value = 124  # was 123
$ws.Range("F13").Value = 330  # was 327
$ws.Range("F14").Value = 2746  # was 2736
$ws.Range("F15").Value = 796  # was 795
$ws.Range("F16").Value = 141  # was 140
$ws.Range("F17").Value = 52  # was 51
$ws.Range("F19").Value = 506  # was 502
$ws.Range("F20").Value = 757  # was 753
$ws.Range("F27").Value = 108  # was 105
$ws.Range("F30").Value = 509  # was 507
$ws.Range("F31").Value = 97  # was 96
$ws.Range("F34").Value = 257  # was 254
$ws.Range("F35").Value = 143  # was 141
$ws.Range("F36").Value = 143  # was 141
$ws.Range("F38").Value = 419  # was 415
$ws.Range("F42").Value = 6  # was 5
$ws.Range("G43").Value = 85  # was 70
$ws.Range("F46").Value = 23  # was 20
$ws.Range("F50").Value = 121  # was 120
